$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H48").Value = 22315
$ws.Range("J48").Value = 22315
$ws.Range("L48").Value = 66945
$ws.Range("N48").Value = -67529
$ws.Range("H56").Value = 22315
$ws.Range("J56").Value = 22315
$ws.Range("L56").Value = 66945
$ws.Range("N56").Value = -68013
$ws.Range("H107").Value = 923.1177
$ws.Range("J107").Value = 753
$ws.Range("L107").Value = 753
$ws.Range("N107").Value = -4593
$ws.Range("H112").Value = 11237554
$ws.Range("J112").Value = 1635.779
$ws.Range("L112").Value = 4907.337
$ws.Range("N112").Value = -7123.337
$ws.Range("H113").Value = 10236.818
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 12200.625
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 12200.625
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -18708.625
$ws.Range("H138").Value = 2301.85
$ws.Range("I138").Value = 1187.0834
$ws.Range("J138").Value = 3974
$ws.Range("K138").Value = 3561.2502
$ws.Range("L138").Value = 11922
$ws.Range("M138").Value = 1578.7498
$ws.Range("N138").Value = -22202

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5985.0107
$ws.Range("I32").Value = 3461.8616
$ws.Range("J32").Value = 11451.833
$ws.Range("K32").Value = 3461.8616
$ws.Range("L32").Value = 11451.833
$ws.Range("M32").Value = -3174.8616
$ws.Range("N32").Value = -12025.833
$ws.Range("H55").Value = 32154.445
$ws.Range("J55").Value = 32154.445
$ws.Range("L55").Value = 32154.445
$ws.Range("N55").Value = -32784.445
$ws.Range("H74").Value = 3716.0857
$ws.Range("I74").Value = 3690.5356
$ws.Range("K74").Value = 3690.5356
$ws.Range("M74").Value = -2816.5356
$ws.Range("H77").Value = 3716.0857
$ws.Range("I77").Value = 3690.5356
$ws.Range("K77").Value = 18452.678
$ws.Range("M77").Value = -14084.678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2457.6667
$ws.Range("I31").Value = 983.65717
$ws.Range("J31").Value = 4802.6816
$ws.Range("K31").Value = 983.65717
$ws.Range("L31").Value = 4802.6816
$ws.Range("M31").Value = -688.65717
$ws.Range("N31").Value = -5392.6816
$ws.Range("H34").Value = 2457.6667
$ws.Range("I34").Value = 983.65717
$ws.Range("J34").Value = 4802.6816
$ws.Range("K34").Value = 983.65717
$ws.Range("L34").Value = 4802.6816
$ws.Range("M34").Value = -781.65717
$ws.Range("N34").Value = -5206.6816
$ws.Range("H35").Value = 23470.5
$ws.Range("I35").Value = 2274.3333
$ws.Range("J35").Value = 44666.668
$ws.Range("K35").Value = 2274.3333
$ws.Range("L35").Value = 44666.668
$ws.Range("M35").Value = -1980.3333
$ws.Range("N35").Value = -45254.668
$ws.Range("H86").Value = 2916.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2916.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2916.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5162.5
$ws.Range("H89").Value = 2916.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2916.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 14582.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -25814.5
$ws.Range("H120").Value = 29999
$ws.Range("J120").Value = 29999
$ws.Range("L120").Value = 29999
$ws.Range("N120").Value = -37257
$ws.Range("H132").Value = 2830.46
$ws.Range("I132").Value = 2041.303
$ws.Range("J132").Value = 4362.353
$ws.Range("K132").Value = 6123.909000000001
$ws.Range("L132").Value = 13087.059
$ws.Range("M132").Value = -3593.909000000001
$ws.Range("N132").Value = -18147.059
$ws.Range("H134").Value = 8426.471
$ws.Range("I134").Value = 10305.091
$ws.Range("J134").Value = 4982.3335
$ws.Range("K134").Value = 30915.273
$ws.Range("L134").Value = 14947.0005
$ws.Range("M134").Value = -28380.273
$ws.Range("N134").Value = -20017.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 18984.125
$ws.Range("I11").Value = 17410.428
$ws.Range("J11").Value = 30000
$ws.Range("K11").Value = 52231.284
$ws.Range("L11").Value = 90000
$ws.Range("M11").Value = -52091.284
$ws.Range("N11").Value = -90280
$ws.Range("H33").Value = 165.9
$ws.Range("I33").Value = 169.75
$ws.Range("J33").Value = 163.33333
$ws.Range("K33").Value = 1018.5
$ws.Range("L33").Value = 979.9999799999999
$ws.Range("M33").Value = -735.5
$ws.Range("N33").Value = -1545.99998
$ws.Range("H113").Value = 653.98334
$ws.Range("I113").Value = 586.67444
$ws.Range("J113").Value = 824.2353000000001
$ws.Range("K113").Value = 1760.02332
$ws.Range("L113").Value = 2472.7059
$ws.Range("M113").Value = 409.97668
$ws.Range("N113").Value = -6812.7059
$ws.Range("H123").Value = 3819.8
$ws.Range("I123").Value = 3525
$ws.Range("K123").Value = 10575
$ws.Range("M123").Value = -8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 28999.6
$ws.Range("J4").Value = 28999.6
$ws.Range("L4").Value = 28999.6
$ws.Range("N4").Value = -29223.6
$ws.Range("H80").Value = 35717056
$ws.Range("I80").Value = 83335660
$ws.Range("K80").Value = 83335660
$ws.Range("M80").Value = -83334662
$ws.Range("H83").Value = 35717056
$ws.Range("I83").Value = 83335660
$ws.Range("K83").Value = 416678300
$ws.Range("M83").Value = -416673308
$ws.Range("H113").Value = 1178.4333
$ws.Range("I113").Value = 1244.5333
$ws.Range("J113").Value = 1112.3334
$ws.Range("K113").Value = 1244.5333
$ws.Range("L113").Value = 1112.3334
$ws.Range("M113").Value = 925.4666999999999
$ws.Range("N113").Value = -5452.3334
$ws.Range("H132").Value = 2691.0334
$ws.Range("I132").Value = 1324.2354
$ws.Range("J132").Value = 4478.385
$ws.Range("K132").Value = 3972.7062
$ws.Range("L132").Value = 13435.155
$ws.Range("M132").Value = -1442.7062
$ws.Range("N132").Value = -18495.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4231.1055
$ws.Range("I7").Value = 2044.7273
$ws.Range("K7").Value = 2044.7273
$ws.Range("M7").Value = -1932.7273
$ws.Range("H126").Value = 4231.1055
$ws.Range("I126").Value = 2044.7273
$ws.Range("K126").Value = 6134.1819
$ws.Range("M126").Value = -3664.1819
$ws.Range("H132").Value = 5096.025
$ws.Range("I132").Value = 1899.8889
$ws.Range("J132").Value = 7711.0454
$ws.Range("K132").Value = 5699.6667
$ws.Range("L132").Value = 23133.1362
$ws.Range("M132").Value = -3169.6667
$ws.Range("N132").Value = -28193.1362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2862.0833
$ws.Range("I122").Value = 2042.0416
$ws.Range("K122").Value = 6126.1248
$ws.Range("M122").Value = -3676.1248
